$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: "Jugendamt" paragraph (sender block) — the address block used
# to end with a run of 27 trailing spaces after "Jugendamt"; that run is
# removed so the paragraph now simply reads "Jugendamt".
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Jugendamt                           ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Jugendamt", 2)

# ---------------------------------------------------------------------
# Edit 2: the "Bitte beachten Sie..." paragraph was split across many
# tiny runs (apparently from sloppy copy/paste). Re-typing the text in
# single Find/Replace passes re-merges the fragmented runs back into
# contiguous runs with uniform formatting, fixing "EBEGU-571 Fehler bei
# Begleitschreiben".
# ---------------------------------------------------------------------

# 2a. "Bitte beachten ... ve" + "r" + "änderte persö" + "n" + "liche Verhältnisse"
#     -> one run for the whole lead-in sentence.
$null = $d.Content.Find.Execute(
    "Bitte beachten Sie auch, dass Sie uns veränderte persönliche Verhältnisse",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Bitte beachten Sie auch, dass Sie uns veränderte persönliche Verhältnisse", 2)

# 2b. Make the single space between "Verhältnisse" and "(z. B." superscript
#     (it was already its own run; keep it that way, just flag the format).
$rng = $d.Content
$null = $rng.Find.Execute("Verhältnisse (z. B.", $true, $false, $false, $false, $false,
                           $true, 1, $false, "", 0)
$spaceStart = $rng.Start + 12
$spaceEnd = $spaceStart + 1
$spaceRng = $d.Range($spaceStart, $spaceEnd)
$spaceRng.Font.Superscript = $true

# 2c. "ren. Ausgenommen davon sind Abwesenheiten wegen Kran" + "k" + "heit, Unfall und während "
$null = $d.Content.Find.Execute(
    "ren. Ausgenommen davon sind Abwesenheiten wegen Krankheit, Unfall und während ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "ren. Ausgenommen davon sind Abwesenheiten wegen Krankheit, Unfall und während ", 2)

# 2d. "Dauer des g" + "e" + "setzlichen Mutterschaftsu" + "r" + "laubs"
$null = $d.Content.Find.Execute(
    "Dauer des gesetzlichen Mutterschaftsurlaubs",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Dauer des gesetzlichen Mutterschaftsurlaubs", 2)

$text = $d.Content.Text
$idx = $text.IndexOf("Jugendamt")
Write-Output $text.Substring($idx, 60)
$idx2 = $text.IndexOf("Bitte beachten")
Write-Output $text.Substring($idx2, 600)
